# "Added last minute updates"
#
# Target paragraph 1 (the **ID__...__ID** placeholder line) gets:
#   - a paragraph border (pBdr top/left/bottom/right, 5-twip space, no line)
#   - its left indent bumped from 120 -> 225 twips (6pt -> 11.25pt)
#   - the placeholder text renamed from
#       **ID__AFFARS_5313_topic_10__ID**
#     to
#       **ID__AFFARS_SUBPART_5313_5__ID**
#   - the trailing " " run removed entirely (no longer just merged away)

$d = $word.ActiveDocument

# --- Remove the lone trailing-space run in paragraph 1 -------------------
# Paragraph 1's text is "**ID__AFFARS_5313_topic_10__ID** " (32 chars of
# placeholder + 1 trailing space) followed by the paragraph mark. Delete
# just the space character/run, leaving the placeholder run untouched.
$p1 = $d.Paragraphs(1)
$spaceRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$spaceRange.Delete()

# --- Rename the placeholder ID -------------------------------------------
$d.Content.Find.Execute("**ID__AFFARS_5313_topic_10__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5313_5__ID**", 2)

# --- Paragraph formatting: indent + border -------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
